$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1443.4286
$ws.Range("I28").Value = 1477.5385
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 1477.5385
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = -992.5385000000001
$ws.Range("N28").Value = -1970
$ws.Range("H29").Value = 2227.5
$ws.Range("I29").Value = 750
$ws.Range("J29").Value = 2523
$ws.Range("K29").Value = 2250
$ws.Range("L29").Value = 7569
$ws.Range("M29").Value = -1969
$ws.Range("N29").Value = -8131
$ws.Range("H32").Value = 1975
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1975
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1975
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -2627
$ws.Range("H41").Value = 565.4666999999999
$ws.Range("I41").Value = 168.875
$ws.Range("J41").Value = 1018.7143
$ws.Range("K41").Value = 168.875
$ws.Range("L41").Value = 1018.7143
$ws.Range("M41").Value = 271.125
$ws.Range("N41").Value = -1898.7143
$ws.Range("H54").Value = 3538
$ws.Range("I54").Value = 3538
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 3538
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -3052
$ws.Range("H117").Value = 19914
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 19914
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 19914
$ws.Range("N117").Value = -29092
$ws.Range("H132").Value = 2714.2
$ws.Range("I132").Value = 2754.5652
$ws.Range("J132").Value = 2250
$ws.Range("K132").Value = 8263.695599999999
$ws.Range("L132").Value = 6750
$ws.Range("M132").Value = -5733.695599999999
$ws.Range("N132").Value = -11810
$ws.Range("H137").Value = 770.2727
$ws.Range("I137").Value = 694
$ws.Range("J137").Value = 833.8333
$ws.Range("K137").Value = 2082
$ws.Range("L137").Value = 2501.4999
$ws.Range("M137").Value = 468
$ws.Range("N137").Value = -7601.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2285.4583
$ws.Range("I2").Value = 1910.5333
$ws.Range("J2").Value = 2910.3333
$ws.Range("K2").Value = 1910.5333
$ws.Range("L2").Value = 2910.3333
$ws.Range("M2").Value = -1797.5333
$ws.Range("H32").Value = 12581.75
$ws.Range("I32").Value = 8145.3823
$ws.Range("J32").Value = 88000
$ws.Range("K32").Value = 8145.3823
$ws.Range("L32").Value = 88000
$ws.Range("M32").Value = -7858.3823
$ws.Range("H63").Value = 1999.3529
$ws.Range("I63").Value = 1999.3125
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 1999.3125
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -1313.3125
$ws.Range("H66").Value = 1999.3529
$ws.Range("I66").Value = 1999.3125
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 9996.5625
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -6564.5625
$ws.Range("H109").Value = 40520
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 40520
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 40520
$ws.Range("N109").Value = -43294
$ws.Range("H116").Value = 2285.4583
$ws.Range("I116").Value = 1910.5333
$ws.Range("J116").Value = 2910.3333
$ws.Range("K116").Value = 1910.5333
$ws.Range("L116").Value = 2910.3333
$ws.Range("M116").Value = 383.4666999999999
$ws.Range("H138").Value = 72370.75
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 72370.75
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 72370.75
$ws.Range("N138").Value = -82650.75
$ws.Range("H140").Value = 83571.75
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 83571.75
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 83571.75
$ws.Range("N140").Value = -93931.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2285.4583
$ws.Range("I3").Value = 1910.5333
$ws.Range("J3").Value = 2910.3333
$ws.Range("K3").Value = 1910.5333
$ws.Range("L3").Value = 2910.3333
$ws.Range("M3").Value = -1796.5333
$ws.Range("H86").Value = 125004670
$ws.Range("I86").Value = 200003330
$ws.Range("J86").Value = 6900
$ws.Range("K86").Value = 200003330
$ws.Range("L86").Value = 6900
$ws.Range("M86").Value = -200002207
$ws.Range("N86").Value = -9146
$ws.Range("H89").Value = 125004670
$ws.Range("I89").Value = 200003330
$ws.Range("J89").Value = 6900
$ws.Range("K89").Value = 1000016650
$ws.Range("L89").Value = 34500
$ws.Range("M89").Value = -1000011034
$ws.Range("N89").Value = -45732
$ws.Range("H94").Value = 799.5217
$ws.Range("I94").Value = 799.5217
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 799.5217
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -348.5217
$ws.Range("H139").Value = 48785
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 48785
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 48785
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -59065

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1654.58
$ws.Range("I31").Value = 1280.8536
$ws.Range("J31").Value = 3357.111
$ws.Range("K31").Value = 1280.8536
$ws.Range("L31").Value = 3357.111
$ws.Range("M31").Value = -985.8535999999999
$ws.Range("N31").Value = -3947.111
$ws.Range("H34").Value = 1654.58
$ws.Range("I34").Value = 1280.8536
$ws.Range("J34").Value = 3357.111
$ws.Range("K34").Value = 1280.8536
$ws.Range("L34").Value = 3357.111
$ws.Range("M34").Value = -1078.8536
$ws.Range("N34").Value = -3761.111
$ws.Range("H42").Value = 5000
$ws.Range("I42").Value = 5000
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 5000
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -4407
$ws.Range("H107").Value = 330.08334
$ws.Range("I107").Value = 330.08334
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 330.08334
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1589.91666
$ws.Range("N107").ClearContents()
$ws.Range("H138").Value = 68876
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 68876
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 68876
$ws.Range("N138").Value = -79156

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 7572
$ws.Range("I44").Value = 400
$ws.Range("J44").Value = 12951
$ws.Range("K44").Value = 1200
$ws.Range("L44").Value = 38853
$ws.Range("M44").Value = -802
$ws.Range("N44").Value = -39649
$ws.Range("H101").Value = 5833.3335
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 5833.3335
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 17500.0005
$ws.Range("N101").Value = -22368.0005
$ws.Range("H131").Value = 14736451
$ws.Range("I131").Value = 45546610
$ws.Range("J131").Value = 1158.6957
$ws.Range("K131").Value = 136639830
$ws.Range("L131").Value = 3476.0871
$ws.Range("M131").Value = -136634790
$ws.Range("N131").Value = -13556.0871
$ws.Range("H133").Value = 2744
$ws.Range("I133").Value = 930
$ws.Range("J133").Value = 10000
$ws.Range("K133").Value = 2790
$ws.Range("L133").Value = 30000
$ws.Range("M133").Value = 2270
$ws.Range("N133").Value = -40120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7666.6665
$ws.Range("I70").Value = 9333.333000000001
$ws.Range("J70").Value = 6000
$ws.Range("K70").Value = 9333.333000000001
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -9063.333000000001
$ws.Range("H73").Value = 7666.6665
$ws.Range("I73").Value = 9333.333000000001
$ws.Range("J73").Value = 6000
$ws.Range("K73").Value = 9333.333000000001
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -8397.333000000001
$ws.Range("H126").Value = 5300.75
$ws.Range("I126").Value = 5903.2
$ws.Range("J126").Value = 4296.6665
$ws.Range("K126").Value = 17709.6
$ws.Range("L126").Value = 12889.9995
$ws.Range("M126").Value = -15239.6
$ws.Range("N126").Value = -17829.9995
$ws.Range("H132").Value = 2922
$ws.Range("I132").Value = 2732.4194
$ws.Range("J132").Value = 3509.7
$ws.Range("K132").Value = 8197.2582
$ws.Range("L132").Value = 10529.1
$ws.Range("M132").Value = -5667.2582
$ws.Range("N132").Value = -15589.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2019.3
$ws.Range("I132").Value = 1354.05
$ws.Range("J132").Value = 3349.8
$ws.Range("K132").Value = 4062.15
$ws.Range("L132").Value = 10049.4
$ws.Range("M132").Value = -1532.15
$ws.Range("N132").Value = -15109.4
$ws.Range("H136").Value = 3739.5588
$ws.Range("I136").Value = 1782.5
$ws.Range("J136").Value = 10100
$ws.Range("K136").Value = 5347.5
$ws.Range("L136").Value = 30300
$ws.Range("M136").Value = -2797.5
$ws.Range("N136").Value = -35400
$ws.Range("H138").Value = 39543
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 39543
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 39543
$ws.Range("N138").Value = -49823

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1002.0417
$ws.Range("I132").Value = 1002.0417
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3006.1251
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -476.1251000000002
$ws.Range("N132").ClearContents()
$ws.Range("H138").Value = 74207.25
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 74207.25
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 74207.25
$ws.Range("N138").Value = -84487.25
$ws.Range("H139").Value = 71715
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 71715
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 71715
$ws.Range("N139").Value = -81995
